$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Daily market-data refresh (coinranking.com symbol list).
# Cells in columns D (Price) and E (Volume 1h %) are numeric-looking text
# (the sheet stores everything as text, e.g. "2.36%"), so each is put into
# Text format before the new value is written - otherwise Excel would
# auto-convert "300.98" / "2.36%" into a real Number/Percentage value.
$updates = @(
    @{ Cell = 'D2'; Value = '300.98' }
    @{ Cell = 'E2'; Value = '2.36%' }
    @{ Cell = 'D3'; Value = '32.23' }
    @{ Cell = 'E3'; Value = '3.12%' }
    @{ Cell = 'D4'; Value = '4.985' }
    @{ Cell = 'E4'; Value = '0.70%' }
    @{ Cell = 'D5'; Value = '0.07745' }
    @{ Cell = 'E5'; Value = '5.68%' }
    @{ Cell = 'D6'; Value = '2.326' }
    @{ Cell = 'E6'; Value = '2.65%' }
    @{ Cell = 'D7'; Value = '7.969' }
    @{ Cell = 'E7'; Value = '3.09%' }
    @{ Cell = 'B8'; Value = 'MXToken' }
    @{ Cell = 'C8'; Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx' }
    @{ Cell = 'D8'; Value = '0.9254' }
    @{ Cell = 'E8'; Value = '1.86%' }
    @{ Cell = 'B9'; Value = 'LiechtensteinCryptoassetsExchange' }
    @{ Cell = 'C9'; Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx' }
    @{ Cell = 'D9'; Value = '0.1005' }
    @{ Cell = 'E9'; Value = '24.38%' }
    @{ Cell = 'B10'; Value = 'WazirX' }
    @{ Cell = 'C10'; Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx' }
    @{ Cell = 'D10'; Value = '0.1759' }
    @{ Cell = 'E10'; Value = '4.35%' }
    @{ Cell = 'B11'; Value = 'MandalaExchangeToken' }
    @{ Cell = 'C11'; Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx' }
    @{ Cell = 'D11'; Value = '0.08468' }
    @{ Cell = 'E11'; Value = '3.90%' }
    @{ Cell = 'B12'; Value = 'BitrueCoin' }
    @{ Cell = 'C12'; Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr' }
    @{ Cell = 'D12'; Value = '0.03292' }
    @{ Cell = 'E12'; Value = '5.85%' }
    @{ Cell = 'B13'; Value = 'BitMartToken' }
    @{ Cell = 'C13'; Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx' }
    @{ Cell = 'D13'; Value = '0.09886' }
    @{ Cell = 'E13'; Value = '-1.89%' }
    @{ Cell = 'B14'; Value = 'BitForexToken' }
    @{ Cell = 'C14'; Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf' }
    @{ Cell = 'D14'; Value = '0.001478' }
    @{ Cell = 'E14'; Value = '-4.67%' }
    @{ Cell = 'B15'; Value = 'TigerCash' }
    @{ Cell = 'C15'; Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch' }
    @{ Cell = 'D15'; Value = '0.005685' }
    @{ Cell = 'E15'; Value = '-0.52%' }
    @{ Cell = 'B16'; Value = 'LEO' }
    @{ Cell = 'C16'; Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo' }
    @{ Cell = 'D16'; Value = '3.510' }
    @{ Cell = 'E16'; Value = '0.81%' }
    @{ Cell = 'B17'; Value = 'GateToken' }
    @{ Cell = 'C17'; Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt' }
    @{ Cell = 'D17'; Value = '3.831' }
    @{ Cell = 'E17'; Value = '2.17%' }
    @{ Cell = 'D18'; Value = '2.193' }
    @{ Cell = 'E18'; Value = '5.56%' }
    @{ Cell = 'D19'; Value = '0.3366' }
    @{ Cell = 'E19'; Value = '1.14%' }
    @{ Cell = 'D20'; Value = '0.1338' }
    @{ Cell = 'E20'; Value = '3.58%' }
    @{ Cell = 'D21'; Value = '4.363' }
    @{ Cell = 'E21'; Value = '9.52%' }
    @{ Cell = 'D23'; Value = '0.04560' }
    @{ Cell = 'E23'; Value = '0.07%' }
    @{ Cell = 'E24'; Value = '0.45%' }
    @{ Cell = 'E25'; Value = '0.61%' }
    @{ Cell = 'E26'; Value = '-0.76%' }
    @{ Cell = 'D27'; Value = '0.0003377' }
    @{ Cell = 'E27'; Value = '-0.74%' }
    @{ Cell = 'D39'; Value = '0.01700' }
    @{ Cell = 'E39'; Value = '6.18%' }
    @{ Cell = 'D40'; Value = '0.04717' }
    @{ Cell = 'E40'; Value = '6.17%' }
    @{ Cell = 'D41'; Value = '0.007733' }
    @{ Cell = 'E41'; Value = '5.07%' }
    @{ Cell = 'D42'; Value = '0.009750' }
    @{ Cell = 'E42'; Value = '12.55%' }
    @{ Cell = 'D43'; Value = '0.1391' }
    @{ Cell = 'E43'; Value = '5.02%' }
    @{ Cell = 'D44'; Value = '0.002096' }
    @{ Cell = 'E44'; Value = '8.20%' }
    @{ Cell = 'D45'; Value = '0.009668' }
    @{ Cell = 'E45'; Value = '1.55%' }
    @{ Cell = 'D46'; Value = '0.00006084' }
    @{ Cell = 'E46'; Value = '2.27%' }
    @{ Cell = 'E47'; Value = '-0.72%' }
    @{ Cell = 'D48'; Value = '2.794' }
    @{ Cell = 'E48'; Value = '24.68%' }
    @{ Cell = 'D49'; Value = '0.001990' }
    @{ Cell = 'E49'; Value = '-31.35%' }
    @{ Cell = 'D50'; Value = '0.00002090' }
    @{ Cell = 'E50'; Value = '-0.72%' }
    @{ Cell = 'D51'; Value = '0.0001990' }
    @{ Cell = 'E51'; Value = '-0.72%' }
)

foreach ($u in $updates) {
    $range = $ws.Range($u.Cell)
    if ($u.Cell[0] -eq "D" -or $u.Cell[0] -eq "E") {
        $range.NumberFormat = "@"
    }
    $range.Value = $u.Value
}
